$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 41551
$ws.Range("B5").Value = 3

# Row 6
$ws.Range("A6").Value = 41552
$ws.Range("B6").Value = 5

# Match the date number format already used in column A (e.g. A2) for the new cells
$ws.Range("A2").Copy()
$ws.Range("A5:A6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Update selection to match diff
$ws.Range("B7").Select()
